$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-22 01:29:44"

# 1) Insert a new blank row at position 11 so the existing "CRM" (row 11) and
#    "CSV" (row 12) rows shift down to 12 and 13, making room for the new
#    MySQL/MariaDB listing.
$ws.Rows("11:11").Insert()

# 2) Refresh the "取得日時" (fetched-at) timestamp on every data row (2-13).
$ws.Range("A2:A13").Value = $newTimestamp

# 3) Rows 5 and 6 swap their listing content (title / price / URL), while the
#    priority score ends up as 333 on both rows.
$ws.Range("B5").Value = "マッチングアプリのLLMO・AIO対策を継続支援いただける方を募集"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5476284"
$ws.Range("G5").Value = 333

$ws.Range("B6").Value = "【急募】マッチングアプリのLLMO・AIO継続支援をお手伝いください!"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5476280"
$ws.Range("G6").Value = 333

# 4) Fill in the newly inserted row 11 with the MySQL/MariaDB listing.
$ws.Range("B11").Value = "【急募】MySQL/MariaDBを活用したデータベース設計・運用の依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5476347"
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = "◇MySQL"

# 5) The row insert leaves the worksheet's hyperlink relationships pointing at
#    their old (pre-insert) rows, so rebuild the whole set from scratch now
#    that every URL cell (F2:F13) holds its final text.
$ws.Range("F2").Hyperlinks.Delete()

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}

Write-Output "done"
